# Add the 4 new raw-submission rows (101-104) to the collection log sheet
# ("八位序列号收集收集结果yd5", the first/active sheet). Columns are:
#   A = submitter name (text)
#   B = submit time (date/time serial, formatted like the existing rows)
#   C = serial number / hex code (text, even when all-digits)
#   D = QQ number (text, even though it looks numeric)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "yyyy/m/d h:mm:ss;@"

$newRows = @(
    @{ Row = 101; A = "Kakarot.";        B = 45969.5273842593; C = "46278631"; D = "3446480369" },
    @{ Row = 102; A = "譚";               B = 45969.6453819445; C = "bfa75387"; D = "2936566799" },
    @{ Row = 103; A = "　  萌虎出没";     B = 45969.8944328704; C = "b9cab845"; D = "2980931837" },
    @{ Row = 104; A = "iiixxxx-7.";      B = 45969.9042361111; C = "53dfa103"; D = "1427967793" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A

    $ws.Cells.Item($row, 2).NumberFormat = $dateFormat
    $ws.Cells.Item($row, 2).Value = $r.B

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $r.C

    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $r.D
}
